# Redesigned File Upload Page
# The "back-end" Heading2 used to read "back-endRegister" (a leftover
# duplication from the "Register" list item that follows it). Clean the
# heading text up so it just reads "back-end".

$d = $word.ActiveDocument

$find = $d.Content.Find
$found = $find.Execute(
    "back-endRegister",  # FindText
    $true,               # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "back-end",           # ReplaceWith
    2                     # Replace (wdReplaceAll)
)

Write-Output "Replaced 'back-endRegister' -> 'back-end': $found"
